$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update maxHealth (column G) values for a handful of enemies
$ws.Range("G9").Value  = 30
$ws.Range("G10").Value = 50
$ws.Range("G11").Value = 50
$ws.Range("G13").Value = 50
$ws.Range("G14").Value = 30
$ws.Range("G15").Value = 30
$ws.Range("G16").Value = 100

# Widen column H (attackRange) to fit the new numbers
$ws.Columns.Item(8).ColumnWidth = 15.86

# Leave the selection on G6, matching the saved cursor position
$ws.Range("G6").Select()
